# Extend the data table by one more year column (2021 -> column O),
# mirroring the formatting already used in column N (2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) from column N into the new column O first,
# row by row, so the new cells inherit the same borders/number formats.
$ws.Range("N3").Copy() | Out-Null
$ws.Range("O3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("N4:N14").Copy() | Out-Null
$ws.Range("O4:O14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# Row 4: new year header
$ws.Range("O4").Value = 2021

# Rows 5-13: new 2021 data values (percent-of-norm figures)
$ws.Range("O5").Value = 70.636215334420882
$ws.Range("O6").Value = 107.1
$ws.Range("O7").Value = 55.452054794520542
$ws.Range("O8").Value = 84.375
$ws.Range("O9").Value = 120.48192771084337
$ws.Range("O10").Value = 109.53346855983774
$ws.Range("O11").Value = 147.7690288713911
$ws.Range("O12").Value = 25.545675020210183
$ws.Range("O13").Value = 82.457854874175425

# Row 14: new 2021 data value
$ws.Range("O14").Value = 15.384615384615385
